$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# 2 more sites entered: F26:F30 flip from "not entered" to "entered"
$ws.Range("F26:F30").Value = "entered"

# Leave the selection where editing finished, on F30
$ws.Range("F30").Select()
